# Added few more test cases
#
# Inserts 3 new rows (for "Space Seperated Input" test cases on the
# No. of Nodes / epsilon / mu fields) right after the existing
# "Enter No. of Nodes field" row, pushing every row below it down by
# three. All of the existing C/D/B content further down the sheet
# naturally shifts along with the inserted rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing row 49 is "Enter No. of Nodes field" / "should accept an
# integer in range 1 to 1000" - the new rows belong right after it, and
# everything from the old row 50 onward (ending at old row 117) needs to
# shift down by three rows.
$ws.Rows("50:52").Insert()

# Row 50 - "Space Seperated Input in No. of Nodes field"
# (write the Description cell first so the new shared strings land in
# the same order as the authored workbook)
$ws.Range("D50").Value = "Should raise an error message"
$ws.Range("C50").Value = "Space Seperated Input in No. of Nodes field"
$ws.Range("E50").Value = "Pass"

# Row 51 - "Space Seperated Input in ε field"
$ws.Range("C51").Value = "Space Seperated Input in ε field"
$ws.Range("D51").Value = "Should raise an error message"
$ws.Range("E51").Value = "Pass"

# Row 52 - "Space Seperated Input in µ field"
$ws.Range("C52").Value = "Space Seperated Input in µ field"
$ws.Range("D52").Value = "Should raise an error message"
$ws.Range("E52").Value = "Pass"

# Match the author's final selection (bottom-right frozen pane active
# cell lands on E51 after the edit).
$ws.Range("E51").Select()
